$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove AutoFilter and the _FilterDatabase defined name ---
$ws.AutoFilterMode = $false
for ($i = $wb.Names.Count; $i -ge 1; $i--) {
    $wb.Names.Item($i).Delete()
}

# --- Delete row 10 (dataset now has one fewer data row: O10 -> O9) ---
$ws.Rows.Item(10).Delete()

# --- Update the description cell (row 1) ---
$ws.Range("A1").Value = "Description unknown, completed 06/15/2023 05:53:58 EDT, by WPJTOWN1.The search returned: 7 events."

# --- Row 3 ---
$ws.Range("A3").Value = "BN"
$ws.Range("B3").Value = 471547
$ws.Range("C3").Value = "Not authorized to view shipment"
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L3").Value = "Not authorized to view shipment"
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("O3").Value = "BN471547"

# --- Row 4 ---
$ws.Range("A4").Value = "CRDX"
$ws.Range("B4").Value = 15008
$ws.Range("C4").Value = "DENVER"
$ws.Range("D4").Value = "CO"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 14
$ws.Range("G4").Value = 303
$ws.Range("H4").Value = "Arrive In-Transit"
$ws.Range("I4").Value = "HKCKDE"
$ws.Range("J4").Value = "LOVELAND"
$ws.Range("K4").Value = "CO"
$ws.Range("L4").Value = 286650
$ws.Range("M4").Value = 68700
$ws.Range("N4").Value = 217950
$ws.Range("O4").Value = "CRDX15008"

# --- Row 5 ---
$ws.Range("A5").Value = "BNSF"
$ws.Range("B5").Value = 468933
$ws.Range("C5").Value = "HOLCOMB"
$ws.Range("D5").Value = "KS"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 13
$ws.Range("G5").Value = 702
$ws.Range("H5").Value = "Departure"
$ws.Range("I5").Value = "HKCKDE"
$ws.Range("J5").Value = "LOVELAND"
$ws.Range("K5").Value = "CO"
$ws.Range("L5").Value = 234960
$ws.Range("M5").Value = 63600
$ws.Range("N5").Value = 171360
$ws.Range("O5").Value = "BNSF468933"

# --- Row 6 ---
$ws.Range("A6").Value = "CRDX"
$ws.Range("B6").Value = 15003
$ws.Range("C6").Value = "JOHNSTOWN"
$ws.Range("D6").Value = "CO"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 12
$ws.Range("G6").Value = 1304
$ws.Range("H6").Value = "Placed Actual"
$ws.Range("I6").ClearContents()
$ws.Range("J6").Value = "LOVELAND"
$ws.Range("K6").Value = "CO"
$ws.Range("L6").Value = 286450
$ws.Range("M6").Value = 68400
$ws.Range("N6").Value = 218050
$ws.Range("O6").Value = "CRDX15003"

# --- Row 7 ---
$ws.Range("A7").Value = "HRTX"
$ws.Range("B7").Value = 541059
$ws.Range("C7").Value = "LITTLETON"
$ws.Range("D7").Value = "CO"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = 101
$ws.Range("H7").Value = "Departure"
$ws.Range("I7").Value = "HKCKDE"
$ws.Range("J7").Value = "LOVELAND"
$ws.Range("K7").Value = "CO"
$ws.Range("L7").Value = 261250
$ws.Range("M7").Value = 64200
$ws.Range("N7").Value = 197050
$ws.Range("O7").Value = "HRTX541059"

# --- Row 8 ---
$ws.Range("A8").Value = "CRDX"
$ws.Range("B8").Value = 15803
$ws.Range("C8").Value = "LOVELAND"
$ws.Range("D8").Value = "CO"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 12
$ws.Range("G8").Value = 1045
$ws.Range("H8").Value = "Junction Received"
$ws.Range("I8").Value = "BNSF"
$ws.Range("J8").Value = "LOVELAND"
$ws.Range("K8").Value = "CO"
$ws.Range("L8").Value = 284700
$ws.Range("M8").Value = 66900
$ws.Range("N8").Value = 217800
$ws.Range("O8").Value = "CRDX15803"

# --- Row 9 ---
$ws.Range("A9").Value = "HRTX"
$ws.Range("B9").Value = 541048
$ws.Range("C9").Value = "LOVELAND"
$ws.Range("D9").Value = "CO"
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 12
$ws.Range("G9").Value = 1045
$ws.Range("H9").Value = "Junction Received"
$ws.Range("I9").Value = "BNSF"
$ws.Range("J9").Value = "LOVELAND"
$ws.Range("K9").Value = "CO"
$ws.Range("L9").Value = 202800
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 202800
$ws.Range("O9").Value = "HRTX541048"

# --- Selection matches new extent ---
$ws.Range("O3:O9").Select()
